# Add viewers annotations to the emfxcel example (studentscourses.xlsx)
#
# - Courses sheet gets a new enrollment row (Gábor Bergmann also takes the
#   "Model-driven software engineering" course, identifier UZ0).
# - Remembered selections on both sheets move on.
# - Courses sheet page setup is touched (portrait / paper size 9).
# - A hyperlink was added to the workbook at some point (and later removed),
#   leaving the built-in "Hyperlink"/"Followed Hyperlink" cell styles behind.

$wb = $excel.ActiveWorkbook
$students = $wb.Worksheets.Item(1)
$courses  = $wb.Worksheets.Item(2)

# --- New data: Gábor Bergmann also enrolls in "Model-driven software engineering" ---
$courses.Range("A8").Value = "Model-driven software engineering"
$courses.Range("B8").Value = "Dániel Varró"
$courses.Range("C8").Value = "UZ0"

# --- Page setup touched on the Courses sheet ---
$coursesPageSetup = $courses.PageSetup
$coursesPageSetup.PaperSize = 9
$coursesPageSetup.Orientation = 1

# --- Leave behind the Hyperlink / Followed Hyperlink style bookkeeping that
#     Excel creates as soon as a hyperlink is inserted anywhere in the
#     workbook, even after the link itself is later removed again. ---
$scratchLink = $students.Range("Z100")
$scratchFollowed = $students.Range("Z101")
$students.Hyperlinks.Add($scratchLink, "http://www.inf.mit.bme.hu/", "", "", "link") | Out-Null
$scratchFollowed.Style = "Followed Hyperlink"

$students.Hyperlinks.Delete()
$scratchLink.ClearFormats() | Out-Null
$scratchLink.ClearContents() | Out-Null
$scratchFollowed.ClearFormats() | Out-Null
$scratchFollowed.ClearContents() | Out-Null

# --- Remembered selections: Courses sheet remembers C9, Students sheet
#     (the active tab) remembers A7. ---
$courses.Range("C9").Select() | Out-Null
$students.Range("A7").Select() | Out-Null
